$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Locate the "Tercera reunion" task-assignment paragraphs up-front, via
# a stable chain of Next() calls, so later edits (which change paragraph
# count / text) cannot shift which paragraph we are pointing at.
# ---------------------------------------------------------------------
$pCristina = $d.Paragraphs(17)   # "Cristina:"
$pJesus    = $pCristina.Next()   # "Jesus:"
$pMiguel   = $pJesus.Next()      # "Miguel:"
$pPablo    = $pMiguel.Next()     # "Pablo: Creacion y configuracion de beans."
$pThalia   = $pPablo.Next()      # "Thalia:" (standalone paragraph to remove)

# --- sanity checks (values only go to the log stream) -----------------
Write-Host ("Cristina before: [" + $pCristina.Range.Text + "]")
Write-Host ("Jesus before: [" + $pJesus.Range.Text + "]")
Write-Host ("Miguel before: [" + $pMiguel.Range.Text + "]")
Write-Host ("Pablo before: [" + $pPablo.Range.Text + "]")
Write-Host ("Thalia before: [" + $pThalia.Range.Text + "]")

# ---------------------------------------------------------------------
# 1) "Cristina:" -> "Cristina y Thalia: creacion de archivos DAO."
# ---------------------------------------------------------------------
$rCristina = $pCristina.Range
$rCristina.MoveEnd(1, -1)              # exclude the paragraph mark
$rCristina.Text = "Cristina y Thalia: creación de archivos DAO."

# ---------------------------------------------------------------------
# 2) "Jesus:" -> "Jesus: creacion de Base de Datos y sus correspondientes tablas."
# ---------------------------------------------------------------------
$rJesus = $pJesus.Range
$rJesus.MoveEnd(1, -1)
$rJesus.InsertAfter(" creación de Base de Datos y sus correspondientes tablas.")

# ---------------------------------------------------------------------
# 3) "Miguel:" -> "Miguel: desarrollo de la capa Negocio." + _GoBack bookmark
#    moved to sit right after this new text (it used to sit after the
#    "Puesta en comun..." paragraph near the end of the document).
#
#    Bookmarks.Add placed at a position that is the very last offset of
#    a run (touching the following element) lands at the wrong spot in
#    this host, so a one-character placeholder is appended, the bookmark
#    is anchored just *before* it (a true mid-run position), and then
#    the placeholder is deleted - the now-collapsed bookmark stays put.
# ---------------------------------------------------------------------
$rMiguel = $pMiguel.Range
$rMiguel.MoveEnd(1, -1)
$rMiguel.InsertAfter(" desarrollo de la capa Negocio.~")

$rMiguelFull = $pMiguel.Range
$rMiguelFull.MoveEnd(1, -1)
$bmAnchor = $d.Range($rMiguelFull.End - 1, $rMiguelFull.End - 1)
$d.Bookmarks.Add("_GoBack", $bmAnchor)

$rMiguelFull2 = $pMiguel.Range
$rMiguelFull2.MoveEnd(1, -1)
$rPlaceholder = $d.Range($rMiguelFull2.End - 1, $rMiguelFull2.End)
$rPlaceholder.Delete()

# ---------------------------------------------------------------------
# 4) Remove the old stand-alone "Thalia:" paragraph entirely (text +
#    its own paragraph mark), so Pablo's paragraph is directly followed
#    by the pre-existing blank paragraph again.
# ---------------------------------------------------------------------
$pThaliaNext = $pThalia.Next()
$rRemove = $d.Range($pThalia.Range.Start, $pThaliaNext.Range.Start)
$rRemove.Delete()

# ---------------------------------------------------------------------
# 5) The old _GoBack bookmark that used to sit after "Puesta en comun
#    ..." has already been relocated by the Bookmarks.Add("_GoBack", ..)
#    call above (bookmark names are unique, so re-adding the same name
#    moves it rather than creating a duplicate) - nothing left to do
#    here, just confirm there is exactly one left and it is in the
#    right place.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $bm = $d.Bookmarks("_GoBack")
    Write-Host ("_GoBack now at Start=" + $bm.Start + " End=" + $bm.End)
}

# --- final sanity check ------------------------------------------------
Write-Host ("Cristina after: [" + $pCristina.Range.Text + "]")
Write-Host ("Jesus after: [" + $pJesus.Range.Text + "]")
Write-Host ("Miguel after: [" + $pMiguel.Range.Text + "]")
Write-Host ("Pablo after: [" + $pPablo.Range.Text + "]")
